$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; all rows from 11 down shift down by one.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new record.
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = [DateTime]"2021-09-14"
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 100112035
$ws.Cells.Item(11, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 40
$ws.Cells.Item(11, 11).Value = 27000
$ws.Cells.Item(11, 12).Value = 27000
$ws.Cells.Item(11, 13).Value = 27000
$ws.Cells.Item(11, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(11, 16).Value = 2700
$ws.Cells.Item(11, 17).Value = 10
$ws.Cells.Item(11, 18).Value = "Hortaliza"
